# Generate Report for Handback
# Adds a new handback row (for file 5221071c-7c60-464c-8500-cfd0f62b8856) to each of
# the three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$fileId   = "5221071c-7c60-464c-8500-cfd0f62b8856"
$fileMd   = "$fileId.md"
$statusOk = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

$rowOv = 4
$wsOverview.Hyperlinks.Add($wsOverview.Range("A$rowOv"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$fileMd", "", "", $fileMd)
$wsOverview.Range("B$rowOv").Value = $statusOk
$wsOverview.Range("C$rowOv").Value = $statusOk

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item(2)

$zhXlf    = "$fileId.f460ea1ad2f1bdb731381398a60430a89e351c3c.zh-cn.xlf"
$zhOffDt  = "2016-02-24 09:40:20"
$zhBackDt = "2016-02-24 09:41:07"
$rowZh = 4

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A$rowZh"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$fileMd", "", "", $fileMd)
$wsZhCn.Range("B$rowZh").Value = $statusOk
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C$rowZh"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf", "", "", $zhXlf)
$wsZhCn.Range("D$rowZh").Value = $zhOffDt
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E$rowZh"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/master/e2e/$fileMd", "", "", $fileMd)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F$rowZh"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf", "", "", $zhXlf)
$wsZhCn.Range("G$rowZh").Value = $zhBackDt
$wsZhCn.Range("H$rowZh").Value = "Include"

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item(3)

$deXlf    = "$fileId.f460ea1ad2f1bdb731381398a60430a89e351c3c.de-de.xlf"
$deOffDt  = "2016-02-24 09:40:32"
$deBackDt = "2016-02-24 09:41:29"
$rowDe = 4

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A$rowDe"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$fileMd", "", "", $fileMd)
$wsDeDe.Range("B$rowDe").Value = $statusOk
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C$rowDe"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf", "", "", $deXlf)
$wsDeDe.Range("D$rowDe").Value = $deOffDt
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E$rowDe"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/master/e2e/$fileMd", "", "", $fileMd)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F$rowDe"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf", "", "", $deXlf)
$wsDeDe.Range("G$rowDe").Value = $deBackDt
$wsDeDe.Range("H$rowDe").Value = "Include"
